# Auto-applies the cryptos list update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text formatting on every cell we touch so that Excel does not
# reinterpret numeric-looking strings (e.g. "1.001", "0.06900") as numbers,
# matching the original inlineStr/text storage used in the workbook.
$cellRefs = @("D2","E2","D3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","E11","E13","D14","E14","D15","E15","D16","E16","D17","E17","D19","E19","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","B40","C40","D40","E40","B41","C41","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","E48","E49","D50","E50","D51","E51")
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the new values
$ws.Range("D2").Value = "28.673.74"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.806.25"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "317.36"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5478"
$ws.Range("E7").Value = "  -4.17%  "
$ws.Range("D8").Value = "0.3807"
$ws.Range("E8").Value = "  -1.96%  "
$ws.Range("D9").Value = "0.07525"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "42.38"
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  -2.01%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").Value = "6.166"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").Value = "7.394"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "1.798.77"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "90.18"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D19").Value = "0.06479"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D21").Value = "17.37"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "5.946"
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "28.679.02"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").Value = "2.096"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "160.24"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").Value = "2.004.35"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").Value = "2.357"
$ws.Range("E29").Value = "  -4.04%  "
$ws.Range("D30").Value = "123.38"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").Value = "1.118"
$ws.Range("E31").Value = "  -3.90%  "
$ws.Range("D32").Value = "0.1060"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "5.657"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").Value = "3.687"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "0.06727"
$ws.Range("E35").Value = "  +9.87%  "
$ws.Range("D36").Value = "0.2258"
$ws.Range("E36").Value = "  +2.09%  "
$ws.Range("D37").Value = "0.02303"
$ws.Range("E37").Value = "  -0.77%  "
$ws.Range("D38").Value = "8.647"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").Value = "5.037"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6258"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "11.29"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("D42").Value = "1.204"
$ws.Range("E42").Value = "  +3.41%  "
$ws.Range("D43").Value = "1.439"
$ws.Range("E43").Value = "  +4.33%  "
$ws.Range("D44").Value = "13.34"
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").Value = "0.5868"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("D46").Value = "3.695"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "126.63"
$ws.Range("E47").Value = "  +3.23%  "
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +1.02%  "
$ws.Range("D50").Value = "0.06900"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").Value = "72.34"
$ws.Range("E51").Value = "  -1.00%  "
